$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new row (row 63) with the next day's gold price data,
# following the same pattern as the existing rows.
$ws.Cells.Item(63, 1).Value = "18-11-2025"
$ws.Cells.Item(63, 2).Value = "The price of gold in India today is ₹12,366 per gram for 24 karat gold, ₹11,335 per gram for 22 karat gold and ₹9,274 per gram for 18 karat gold (also called 999 gold)."
